{"js": "// Office.js (Word JavaScript API) script\n// Body is the implementation of: async (context) => { ... }\n\n// Helper: build a minimal OOXML \"flat package\" wrapper so that Range.insertOoxml\n// can be used to replace a Range's content with a specific sequence of <w:r>/<w:proofErr>\n// elements while leaving the *hosting* paragraph's <w:pPr> (list style, numbering, etc.)\n// completely untouched (because we only ever replace an inner text Range, never the\n// paragraph mark itself).\nfunction wrapOoxml(innerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + innerXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// Standard run properties used throughout this document.\nconst RPR = '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>';\n\nfunction run(text, preserveSpace) {\n  const space = preserveSpace ? ' xml:space=\"preserve\"' : '';\n  return '<w:r>' + RPR + '<w:t' + space + '>' + text + '</w:t></w:r>';\n}\n\nfunction gram(text, preserveSpace) {\n  return '<w:proofErr w:type=\"gramStart\"/>' + run(text, preserveSpace) + '<w:proofErr w:type=\"gramEnd\"/>';\n}\n\nfunction spell(text, preserveSpace) {\n  return '<w:proofErr w:type=\"spellStart\"/>' + run(text, preserveSpace) + '<w:proofErr w:type=\"spellEnd\"/>';\n}\n\nasync function replaceParagraphText(body, exactText, innerXmlBuilder) {\n  const results = body.search(exactText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + exactText);\n  }\n  const range = results.items[0];\n  range.insertOoxml(wrapOoxml(innerXmlBuilder()), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\n\n// 1) \"Configure the build path to include Tomcat and the JDBC connector file\"\nawait replaceParagraphText(\n  body,\n  \"Configure the build path to include Tomcat and the JDBC connector file\",\n  () =>\n    run(\"Configure the build path to include Tomcat and the JDBC connector \", true) +\n    gram(\"file\")\n);\nawait context.sync();\n\n// 2) \"Create basic HTML file to include in the WebINF folder\"\nawait replaceParagraphText(\n  body,\n  \"Create basic HTML file to include in the WebINF folder\",\n  () =>\n    run(\"Create basic HTML file to include in the \", true) +\n    spell(\"WebINF\") +\n    run(\" \", true) +\n    gram(\"folder\")\n);\nawait context.sync();\n\n// 3) \"From the original Class Diagram, the Login and Registration classes will be Servlets\"\nawait replaceParagraphText(\n  body,\n  \"From the original Class Diagram, the Login and Registration classes will be Servlets\",\n  () =>\n    run(\"From the original Class Diagram, the Login and Registration classes will be \", true) +\n    gram(\"Servlets\")\n);\nawait context.sync();\n\n// 4) \"Create the database\"\nawait replaceParagraphText(\n  body,\n  \"Create the database\",\n  () => run(\"Create the \", true) + gram(\"database\")\n);\nawait context.sync();\n\n// 5) \" Workbench\" (second run of \"Make it manually in MySQL Workbench\" paragraph)\nawait replaceParagraphText(\n  body,\n  \" Workbench\",\n  () => run(\" \", true) + gram(\"Workbench\")\n);\nawait context.sync();\n\n// 6) \"In order to collaborate, we are using Github (through the native Eclipse features)\"\nawait replaceParagraphText(\n  body,\n  \"In order to collaborate, we are using Github (through the native Eclipse features)\",\n  () =>\n    gram(\"In order to\") +\n    run(\" collaborate, we are \", true) +\n    run(\"using \", true) +\n    spell(\"Github\") +\n    run(\" (through the native Eclipse features)\", true)\n);\nawait context.sync();\n\n// 7) \"Learning how to get Git to work in Eclipse was difficult for all of us\"\nawait replaceParagraphText(\n  body,\n  \"Learning how to get Git to work in Eclipse was difficult for all of us\",\n  () =>\n    run(\"Learning how to get Git to work in Eclipse was difficult for all of \", true) +\n    gram(\"us\")\n);\nawait context.sync();\n\n// 8) Delete the seven paragraphs that followed (no longer part of the report):\nconst paragraphsToDelete = [\n  \"Started by coding together, after 2 meetings decided to break up the workload somewhat\",\n  \"Servlets for Erica (Registration and Login)\",\n  \"DAO for Kara (who already did the HTML forms)\",\n  \"Beans for Caitlin (User and Guardian classes, plus updating report)\",\n  \"Continue to have issues getting Git to work\",\n  \"Issues seem to stem from trying to commit changes before pulling the repository properly\",\n  \"As a workaround, we\\u2019ll try making separate branches\"\n];\n\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"text\");\nawait context.sync();\n\nfor (const target of paragraphsToDelete) {\n  let found = false;\n  for (let i = 0; i < allParagraphs.items.length; i++) {\n    if (allParagraphs.items[i].text === target) {\n      allParagraphs.items[i].delete();\n      found = true;\n      break;\n    }\n  }\n  if (!found) {\n    throw new Error(\"Could not find paragraph to delete: \" + target);\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# Helper: locate an exact phrase anywhere in the document, then replace that\n# exact range's content with freshly supplied OOXML runs (no paragraph mark\n# included), so the hosting paragraph's <w:pPr> (numbering / style) is left\n# completely untouched.\n# ---------------------------------------------------------------------------\nfunction Replace-TextWithOoxml($SearchText, $InnerXml) {\n    $findRange = $d.Content\n    $found = $findRange.Find.Execute($SearchText)\n    if (-not $found) {\n        throw \"Could not find text: $SearchText\"\n    }\n\n    # Build a brand-new Range object from the located Start/End; calling\n    # InsertXML directly on the Range mutated by Find.Execute does not\n    # perform a clean replace, so we re-wrap the bounds first.\n    $target = $d.Range($findRange.Start, $findRange.End)\n\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' + $InnerXml + '</w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    $target.InsertXML($xml)\n}\n\n$RPR = '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>'\n\nfunction Make-Run($Text, $PreserveSpace) {\n    $space = \"\"\n    if ($PreserveSpace) { $space = ' xml:space=\"preserve\"' }\n    return \"<w:r>$RPR<w:t$space>$Text</w:t></w:r>\"\n}\n\nfunction Make-Gram($Text, $PreserveSpace) {\n    return '<w:proofErr w:type=\"gramStart\"/>' + (Make-Run $Text $PreserveSpace) + '<w:proofErr w:type=\"gramEnd\"/>'\n}\n\nfunction Make-Spell($Text, $PreserveSpace) {\n    return '<w:proofErr w:type=\"spellStart\"/>' + (Make-Run $Text $PreserveSpace) + '<w:proofErr w:type=\"spellEnd\"/>'\n}\n\n# 1) \"Configure the build path to include Tomcat and the JDBC connector file\"\n$inner = (Make-Run \"Configure the build path to include Tomcat and the JDBC connector \" $true) +\n         (Make-Gram \"file\" $false)\nReplace-TextWithOoxml \"Configure the build path to include Tomcat and the JDBC connector file\" $inner\n\n# 2) \"Create basic HTML file to include in the WebINF folder\"\n$inner = (Make-Run \"Create basic HTML file to include in the \" $true) +\n         (Make-Spell \"WebINF\" $false) +\n         (Make-Run \" \" $true) +\n         (Make-Gram \"folder\" $false)\nReplace-TextWithOoxml \"Create basic HTML file to include in the WebINF folder\" $inner\n\n# 3) \"From the original Class Diagram, the Login and Registration classes will be Servlets\"\n$inner = (Make-Run \"From the original Class Diagram, the Login and Registration classes will be \" $true) +\n         (Make-Gram \"Servlets\" $false)\nReplace-TextWithOoxml \"From the original Class Diagram, the Login and Registration classes will be Servlets\" $inner\n\n# 4) \"Create the database\"\n$inner = (Make-Run \"Create the \" $true) + (Make-Gram \"database\" $false)\nReplace-TextWithOoxml \"Create the database\" $inner\n\n# 5) \" Workbench\" (second run of \"Make it manually in MySQL Workbench\")\n$inner = (Make-Run \" \" $true) + (Make-Gram \"Workbench\" $false)\nReplace-TextWithOoxml \" Workbench\" $inner\n\n# 6) \"In order to collaborate, we are using Github (through the native Eclipse features)\"\n$inner = (Make-Gram \"In order to\" $false) +\n         (Make-Run \" collaborate, we are \" $true) +\n         (Make-Run \"using \" $true) +\n         (Make-Spell \"Github\" $false) +\n         (Make-Run \" (through the native Eclipse features)\" $true)\nReplace-TextWithOoxml \"In order to collaborate, we are using Github (through the native Eclipse features)\" $inner\n\n# 7) \"Learning how to get Git to work in Eclipse was difficult for all of us\"\n$inner = (Make-Run \"Learning how to get Git to work in Eclipse was difficult for all of \" $true) +\n         (Make-Gram \"us\" $false)\nReplace-TextWithOoxml \"Learning how to get Git to work in Eclipse was difficult for all of us\" $inner\n\n# ---------------------------------------------------------------------------\n# 8) Remove the seven paragraphs that used to follow (no longer part of the\n#    report). Match by exact paragraph text (trimming the trailing paragraph\n#    mark) and walk backwards so deleting does not disturb earlier indices.\n# ---------------------------------------------------------------------------\n$paragraphsToDelete = @(\n    \"Started by coding together, after 2 meetings decided to break up the workload somewhat\",\n    \"Servlets for Erica (Registration and Login)\",\n    \"DAO for Kara (who already did the HTML forms)\",\n    \"Beans for Caitlin (User and Guardian classes, plus updating report)\",\n    \"Continue to have issues getting Git to work\",\n    \"Issues seem to stem from trying to commit changes before pulling the repository properly\",\n    (\"As a workaround, we\" + [char]0x2019 + \"ll try making separate branches\")\n)\n\nforeach ($target in $paragraphsToDelete) {\n    $found = $false\n    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n        $p = $d.Paragraphs.Item($i)\n        $ptext = $p.Range.Text.TrimEnd([char]13)\n        if ($ptext -eq $target) {\n            $p.Range.Delete()\n            $found = $true\n            break\n        }\n    }\n    if (-not $found) {\n        throw \"Could not find paragraph to delete: $target\"\n    }\n}\n"}
